$wb = $excel.ActiveWorkbook

# Rename the "SwateTemplateMetadata" sheet to "isa_template"
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "SwateTemplateMetadata") {
        $sheet.Name = "isa_template"
    }
}
